$wb = $excel.ActiveWorkbook

# --- Sheet "rfid_item": replace the single-letter RFID tags with real tag IDs ---
$ws2 = $wb.Worksheets.Item("rfid_item")
$ws2.Range("A2").Value = "C7E671B4"
$ws2.Range("A3").Value = "0A4B997F"
$ws2.Range("A4").Value = "AAAAAAA"
$ws2.Range("A5").Value = "BBBBBBBB"
$ws2.Range("A6").Value = "CCCCCCCC"
$ws2.Range("A7").Value = "DDDDDDDD"

# Auto-fit column A on rfid_item to match the new, wider content
$ws2.Columns.Item(1).AutoFit() | Out-Null

# --- Sheet "device_cart": update the device name in row 2 ---
$ws1 = $wb.Worksheets.Item("device_cart")
$ws1.Range("A2").Value = "zovag"

# --- Restore the selections left by the editor in each sheet ---
$ws1.Range("A2").Select() | Out-Null
$ws2.Range("A3").Select() | Out-Null
